# Applies the LOB1258.xlsx content edit described by the commit diff:
# rewrites the "Objetivos/Programa/Avaliacao" block (rows 10, 13-23) to the
# corrected/reordered text and drops the trailing Requisitos/Bibliografia rows
# (24-25) that no longer exist afterwards (final used range becomes A1:C23).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two trailing rows that disappear in the target layout.
# Deleting from the bottom up keeps the row numbers of everything above stable.
$ws.Rows.Item(25).Delete()
$ws.Rows.Item(24).Delete()

# --- Row 10 ---
$ws.Range("A10").Value = "Objetivos:"
$ws.Range("B10").Value = "4780627 - Ana Lucia Gabas Ferreira"
$ws.Range("C10").Value = "4780627 - Ana Lucia Gabas Ferreira"
$ws.Rows.Item(10).RowHeight = 60

# --- Row 13 ---
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "01/01/2022"
$ws.Range("C13").Value = "01/01/2022"
$ws.Rows.Item(13).RowHeight = 60

# --- Row 14 ---
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Value = "Permanent flow in pressure conduits. Head losses in pipelines. Pumping stations and hydraulic pumps. Permanente flow in free surface conduits, head losses in natural and artificial channels. Specific Energy in channels."
$ws.Range("C14").Value = "Permanent flow in pressure conduits. Head losses in pipelines. Pumping stations and hydraulic pumps. Permanente flow in free surface conduits, head losses in natural and artificial channels. Specific Energy in channels."
$ws.Rows.Item(14).RowHeight = 60

# --- Row 15 ---
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "4780627 - Ana Lucia Gabas Ferreira"
$ws.Range("C15").Value = "4780627 - Ana Lucia Gabas Ferreira"
$ws.Rows.Item(15).RowHeight = 120

# --- Row 16 ---
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Value = "- Hydrostatic,- piezometry,- conservation of mass and amount of movement,- Permanent flows in Pipes,- Resistance to Flow and Load Losses,- Pumps and discharge systems,- Uniform Permanent Flow in Free Flues,- Resistance to flow and head losses,- Regular and natural channels,- Specific Load,- Gradually Variable Permanent Flow,- Calculation of the water line,- Hydraulic boss."
$ws.Range("C16").Value = "- Hydrostatic,- piezometry,- conservation of mass and amount of movement,- Permanent flows in Pipes,- Resistance to Flow and Load Losses,- Pumps and discharge systems,- Uniform Permanent Flow in Free Flues,- Resistance to flow and head losses,- Regular and natural channels,- Specific Load,- Gradually Variable Permanent Flow,- Calculation of the water line,- Hydraulic boss."
$ws.Rows.Item(16).RowHeight = 120

# --- Row 17 ---
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("B17").Clear()
$ws.Range("C17").Clear()
$ws.Rows.Item(17).AutoFit()

# --- Row 18 ---
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "7455355 - Robson da Silva Rocha"
$ws.Range("C18").Value = "7455355 - Robson da Silva Rocha"
$ws.Rows.Item(18).RowHeight = 60

# --- Row 19 ---
$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "Aulas teóricas e práticas, trabalhos de campo e exercícios dirigidos.Avaliação baseada em provas, exercícios e trabalhos práticos e relatórios."
$ws.Range("C19").Value = "Aulas teóricas e práticas, trabalhos de campo e exercícios dirigidos.Avaliação baseada em provas, exercícios e trabalhos práticos e relatórios."
$ws.Rows.Item(19).RowHeight = 60

# --- Row 20 ---
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "Média ponderada das notas atribuídas às provas, exercícios e trabalhos práticos e relatórios."
$ws.Range("C20").Value = "Média ponderada das notas atribuídas às provas, exercícios e trabalhos práticos e relatórios."
$ws.Rows.Item(20).RowHeight = 60

# --- Row 21 ---
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "1 (uma) prova de recuperação (R), sendo considerado aprovado se 0,5(NF + R) >= 5,0."
$ws.Range("C21").Value = "1 (uma) prova de recuperação (R), sendo considerado aprovado se 0,5(NF + R) >= 5,0."
$ws.Rows.Item(21).RowHeight = 120

# --- Row 22 ---
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B22").Clear()
$ws.Range("C22").Clear()
$ws.Rows.Item(22).AutoFit()

# --- Row 23 ---
$ws.Range("A23").Clear()
$ws.Range("B23").Value = "LOB1217 -  Operações Unitárias e Processos  (Requisito fraco)`n"
$ws.Range("C23").Value = "LOB1217 -  Operações Unitárias e Processos  (Requisito fraco)`n"
$ws.Rows.Item(23).RowHeight = 30
